$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F ("język") - this shifts columns G:J left to F:I
$ws.Columns("F").Delete()

# Restore the active cell selection to F1 as per the saved file
$ws.Range("F1").Select()

# The hyperlinks that lived on the "email" column (old column I) need to be
# re-pointed at their new location (now column H) since deleting a column
# does not automatically re-anchor existing Hyperlink objects.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 4; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $ws.Hyperlinks.Add($cell, "mailto:ala@gmail.com", [Type]::Missing, [Type]::Missing, "ala@gmail.com")
}
